$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template source cells that already carry the desired border (thin, all sides)
# and will have their alignment adjusted after the format copy so no new border
# style gets synthesized.
$fmtSrcPlain = $ws.Range("A965:G965")
$fmtSrcFail  = $ws.Range("D965")

# ---- Row 1018 ----
$ws.Range("A1018").Value = "MH"
$ws.Range("B1018").Value = "NAG6426_NGP_P40"
$ws.Range("C1018").Value = "27-Dec-2025 6:03 AM"
$ws.Range("D1018").Value = "FAIL"
$ws.Range("E1018").Value = "1. SCG addition after VoLTE call released"
$ws.Range("F1018").Value = "1. Static VoLTE MO"
$ws.Range("G1018").Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition."
$fmtSrcPlain.Copy()
$ws.Range("A1018:G1018").PasteSpecial(-4122)
$ws.Range("A1018:G1018").HorizontalAlignment = 1
$ws.Range("A1018:G1018").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1018").PasteSpecial(-4122)
$ws.Range("D1018").HorizontalAlignment = 1
$ws.Range("D1018").VerticalAlignment = -4160
$ws.Rows.Item(1018).RowHeight = 87

# ---- Row 1019 ----
$ws.Range("A1019").Value = "TN"
$ws.Range("B1019").Value = "CB1221_CBE_P40"
$ws.Range("C1019").Value = "27-Dec-2025 8:22 AM"
$ws.Range("D1019").Value = "FAIL"
$ws.Range("E1019").Value = "1. Video Streaming  (ms)"
$ws.Range("F1019").Value = "1. Static Yotube Streaming"
$ws.Range("G1019").Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrcPlain.Copy()
$ws.Range("A1019:G1019").PasteSpecial(-4122)
$ws.Range("A1019:G1019").HorizontalAlignment = 1
$ws.Range("A1019:G1019").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1019").PasteSpecial(-4122)
$ws.Range("D1019").HorizontalAlignment = 1
$ws.Range("D1019").VerticalAlignment = -4160
$ws.Rows.Item(1019).RowHeight = 58

# ---- Row 1020 ----
$ws.Range("A1020").Value = "KL"
$ws.Range("B1020").Value = "VENG29_EKM_P40"
$ws.Range("C1020").Value = "26-Dec-2025 9:15 PM"
$ws.Range("D1020").Value = "FAIL"
$ws.Range("E1020").Value = "1. SCG addition after VoLTE call released
2. SgNB Addition time (ms)"
$ws.Range("F1020").Value = "1. Static VoLTE MO
2. Static ATDT"
$ws.Range("G1020").Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.
2. Exclue ATDT Logfile and Create New Sgnb Addition Time Is Very High. It Should Be <150 Ms. To Achieve This, Perform Static Test In Main Lobe And Keep Test Files Downloading In Background. Also, Ensure 4G Serving Cell Belongs To The Same Site. Exclude The Existing Logfile First"
$fmtSrcPlain.Copy()
$ws.Range("A1020:G1020").PasteSpecial(-4122)
$ws.Range("A1020:G1020").HorizontalAlignment = 1
$ws.Range("A1020:G1020").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1020").PasteSpecial(-4122)
$ws.Range("D1020").HorizontalAlignment = 1
$ws.Range("D1020").VerticalAlignment = -4160
$ws.Rows.Item(1020).RowHeight = 145

# ---- Row 1021 ----
$ws.Range("A1021").Value = "PB"
$ws.Range("B1021").Value = "MIR803_AMS_P40"
$ws.Range("C1021").Value = "26-Dec-2025 9:02 PM"
$ws.Range("D1021").Value = "FAIL"
$ws.Range("E1021").Value = "1. SCG addition after VoLTE call released
2. Peak Rank - 5G
3. Peak PUSCH UL Throughput"
$ws.Range("F1021").Value = "1. Static VoLTE MO
2. Static DL
3. Static UL"
$ws.Range("G1021").Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.
2. Peak Rank is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.
3. Peak PUSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PUSCH Throughput in the NR tab."
$fmtSrcPlain.Copy()
$ws.Range("A1021:G1021").PasteSpecial(-4122)
$ws.Range("A1021:G1021").HorizontalAlignment = 1
$ws.Range("A1021:G1021").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1021").PasteSpecial(-4122)
$ws.Range("D1021").HorizontalAlignment = 1
$ws.Range("D1021").VerticalAlignment = -4160
$ws.Rows.Item(1021).RowHeight = 203

# ---- Row 1022 ----
$ws.Range("A1022").Value = "MP"
$ws.Range("B1022").Value = "MPID8973_IND_P40"
$ws.Range("C1022").Value = "26-Dec-2025 9:12 PM"
$ws.Range("D1022").Value = "FAIL"
$ws.Range("E1022").Value = "1. Peak PUSCH UL Throughput"
$ws.Range("F1022").Value = "1. Static UL"
$ws.Range("G1022").Value = "1. Peak PUSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PUSCH Throughput in the NR tab."
$fmtSrcPlain.Copy()
$ws.Range("A1022:G1022").PasteSpecial(-4122)
$ws.Range("A1022:G1022").HorizontalAlignment = 1
$ws.Range("A1022:G1022").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1022").PasteSpecial(-4122)
$ws.Range("D1022").HorizontalAlignment = 1
$ws.Range("D1022").VerticalAlignment = -4160
$ws.Rows.Item(1022).RowHeight = 58

# ---- Row 1023 ----
$ws.Range("A1023").Value = "KL"
$ws.Range("B1023").Value = "PKD175_PKD_P40"
$ws.Range("C1023").Value = "26-Dec-2025 8:06 PM"
$ws.Range("D1023").Value = "FAIL"
$ws.Range("E1023").Value = "1. SCG addition after VoLTE call released"
$ws.Range("F1023").Value = "1. Static VoLTE MO"
$ws.Range("G1023").Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition."
$fmtSrcPlain.Copy()
$ws.Range("A1023:G1023").PasteSpecial(-4122)
$ws.Range("A1023:G1023").HorizontalAlignment = 1
$ws.Range("A1023:G1023").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1023").PasteSpecial(-4122)
$ws.Range("D1023").HorizontalAlignment = 1
$ws.Range("D1023").VerticalAlignment = -4160
$ws.Rows.Item(1023).RowHeight = 87

# ---- Row 1024 ----
$ws.Range("A1024").Value = "KL"
$ws.Range("B1024").Value = "ANCL05_KLM_P40"
$ws.Range("C1024").Value = "26-Dec-2025 7:16 PM"
$ws.Range("D1024").Value = "FAIL"
$ws.Range("E1024").Value = "1. SgNB Addition time (ms)
2. Video Streaming  (ms)"
$ws.Range("F1024").Value = "1. Static ATDT
2. Static Yotube Streaming"
$ws.Range("G1024").Value = "1. Exclue ATDT Logfile and Create New Sgnb Addition Time Is Very High. It Should Be <150 Ms. To Achieve This, Perform Static Test In Main Lobe And Keep Test Files Downloading In Background. Also, Ensure 4G Serving Cell Belongs To The Same Site. Exclude The Existing Logfile First
2. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrcPlain.Copy()
$ws.Range("A1024:G1024").PasteSpecial(-4122)
$ws.Range("A1024:G1024").HorizontalAlignment = 1
$ws.Range("A1024:G1024").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1024").PasteSpecial(-4122)
$ws.Range("D1024").HorizontalAlignment = 1
$ws.Range("D1024").VerticalAlignment = -4160
$ws.Rows.Item(1024).RowHeight = 116

# ---- Row 1025 ----
$ws.Range("A1025").Value = "RJ"
$ws.Range("B1025").Value = "VJKN02_JPR_P40"
$ws.Range("C1025").Value = "26-Dec-2025 5:45 PM"
$ws.Range("D1025").Value = "FAIL"
$ws.Range("E1025").Value = "1. SCG addition after VoLTE call released"
$ws.Range("F1025").Value = "1. Static VoLTE MO"
$ws.Range("G1025").Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition."
$fmtSrcPlain.Copy()
$ws.Range("A1025:G1025").PasteSpecial(-4122)
$ws.Range("A1025:G1025").HorizontalAlignment = 1
$ws.Range("A1025:G1025").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1025").PasteSpecial(-4122)
$ws.Range("D1025").HorizontalAlignment = 1
$ws.Range("D1025").VerticalAlignment = -4160
$ws.Rows.Item(1025).RowHeight = 87

# ---- Row 1026 ----
$ws.Range("A1026").Value = "BH"
$ws.Range("B1026").Value = "BHDAL-25_PAT_P41"
$ws.Range("C1026").Value = "26-Dec-2025 5:42 PM"
$ws.Range("D1026").Value = "FAIL"
$ws.Range("E1026").Value = "1. Reselection Priorities in Respective RAT's
2. RACH setup
3. Ping/Round trip time(ms)
4. Video Streaming"
$ws.Range("F1026").Value = "1. Static All
2. Static All
3. Static Ping
4. Static Yotube Streaming"
$ws.Range("G1026").Value = "1. These parameters are auto-derived and generally do not fail. As the remark is not available, please check with the report provider or ANTS Support Team.
2. If DT Tool is  TEMS Pocket, verify the Static ATDT . The NR RACH Attempts should be equal to NR RACH Success; kindly exclude the logs where NR RACH has failed.
If DT Tool is  AZQ, Validate Static All and ensure NR RACH Attempts match NR RACH Success. Please exclude the logs with NR RACH failures and redo the test accordingly.
3. Ping is not meeting the acceptance criteria. The average ping value across all logfiles should be less than 50 ms. Kindly exclude the logfile where the average value exceeds 50 ms and redo the test.”
4. While performing the YouTube test for both sectors, please ensure that the video is playing successfully in the script before saving the log file."
$fmtSrcPlain.Copy()
$ws.Range("A1026:G1026").PasteSpecial(-4122)
$ws.Range("A1026:G1026").HorizontalAlignment = 1
$ws.Range("A1026:G1026").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1026").PasteSpecial(-4122)
$ws.Range("D1026").HorizontalAlignment = 1
$ws.Range("D1026").VerticalAlignment = -4160
$ws.Rows.Item(1026).RowHeight = 203

# ---- Row 1027 ----
$ws.Range("A1027").Value = "BH"
$ws.Range("B1027").Value = "BHHAJ-103_PAT_P41"
$ws.Range("C1027").Value = "26-Dec-2025 5:18 PM"
$ws.Range("D1027").Value = "FAIL"
$ws.Range("E1027").Value = "1. Peak PDSCH DL Throughput
2. Peak Rank - 4G
3. SgNB Addition time (ms)
4. UE Steering (Idle) : Non anchor/anchor to preferred anchor
5. Video Streaming"
$ws.Range("F1027").Value = "1. Static DL
2. Static DL
3. Static ATDT
4. Static Idle
5. Static Yotube Streaming"
$ws.Range("G1027").Value = "1. Peak PDSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PDSCH Throughput in the NR tab.
2. Peak Rank is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.
3. Exclue ATDT Logfile and Create New Sgnb Addition Time Is Very High. It Should Be <150 Ms. To Achieve This, Perform Static Test In Main Lobe And Keep Test Files Downloading In Background. Also, Ensure 4G Serving Cell Belongs To The Same Site. Exclude The Existing Logfile First
4. For sites with NOKIA OEM, validate using Drive Idle, and for other OEMs, validate using Static Idle. In both Drive and Static Idle, the UE should latch from NR to LTE and from LTE to NR. In LTE, the UE should latch on the band that corresponds to the configured anchor layer.
5. While performing the YouTube test for both sectors, please ensure that the video is playing successfully in the script before saving the log file."
$fmtSrcPlain.Copy()
$ws.Range("A1027:G1027").PasteSpecial(-4122)
$ws.Range("A1027:G1027").HorizontalAlignment = 1
$ws.Range("A1027:G1027").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1027").PasteSpecial(-4122)
$ws.Range("D1027").HorizontalAlignment = 1
$ws.Range("D1027").VerticalAlignment = -4160
$ws.Rows.Item(1027).RowHeight = 261

# ---- Row 1028 ----
$ws.Range("A1028").Value = "TN"
$ws.Range("B1028").Value = "PAPNP3_VLP_P40"
$ws.Range("C1028").Value = "26-Dec-2025 2:13 PM"
$ws.Range("D1028").Value = "FAIL"
$ws.Range("E1028").Value = "1. Video Streaming  (ms)"
$ws.Range("F1028").Value = "1. Static Yotube Streaming"
$ws.Range("G1028").Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrcPlain.Copy()
$ws.Range("A1028:G1028").PasteSpecial(-4122)
$ws.Range("A1028:G1028").HorizontalAlignment = 1
$ws.Range("A1028:G1028").VerticalAlignment = -4160
$fmtSrcFail.Copy()
$ws.Range("D1028").PasteSpecial(-4122)
$ws.Range("D1028").HorizontalAlignment = 1
$ws.Range("D1028").VerticalAlignment = -4160
$ws.Rows.Item(1028).RowHeight = 58

$excel.CutCopyMode = $false
$ws.Range("A1018:G1028").Select()
